$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New rows describing the Lancaster Ave daily traffic volume source
$ws.Range("A18").Value = "Daily Traffic Volume"

$ws.Range("A19").Value = "Lancaster Ave"
$ws.Range("B19").Value = 15000

$trafficUrl = "https://gis.penndot.pa.gov/BPR_PDF_FILES/MAPS/Traffic/Traffic_Volume/Statewide/Statewide_2022_tv.pdf"
$ws.Range("D19").Value = $trafficUrl
$ws.Hyperlinks.Add($ws.Range("D19"), $trafficUrl) | Out-Null
$ws.Range("D19").Style = "Hyperlink"

# Widen column A (school/category names) and D (source links) to fit the new content
$ws.Columns.Item(1).ColumnWidth = 27.3
$ws.Columns.Item(4).ColumnWidth = 89.8

# Selection moves past the new rows
$ws.Range("A20").Select()
